# Applies the coin-table refresh for cryptos.xlsx (prices/volumes updated,
# Mantle <-> FirstDigitalUSD rows swapped) as of the Aug 14 2024 GitHub Action run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.844.22"
$ws.Range("E2").Value = "  -2.53%  "

$ws.Range("D3").Value = "2.656.99"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'523.08"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").Value = "'143.94"
$ws.Range("E6").Value = "  -2.21%  "

$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("E9").Value = "  +6.61%  "

$ws.Range("D10").Value = "'0.102"
$ws.Range("E10").Value = "  -3.84%  "

$ws.Range("D11").Value = "'0.335"
$ws.Range("E11").Value = "  -1.46%  "

$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").Value = "3.127.95"
$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("D14").Value = "58.858.48"
$ws.Range("E14").Value = "  -2.54%  "

$ws.Range("D15").Value = "'20.99"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").Value = "2.661.78"
$ws.Range("E17").Value = "  -6.75%  "

$ws.Range("D18").Value = "'338.79"
$ws.Range("E18").Value = "  -3.63%  "

$ws.Range("E19").Value = "  -3.61%  "

$ws.Range("D20").Value = "'10.38"
$ws.Range("E20").Value = "  -2.08%  "

$ws.Range("D21").Value = "'6.35"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").Value = "'64.43"
$ws.Range("E23").Value = "  +2.03%  "

$ws.Range("D24").Value = "'0.420"
$ws.Range("E24").Value = "  -0.84%  "

$ws.Range("D25").Value = "'0.166"
$ws.Range("E25").Value = "  -1.60%  "

$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.35%  "

$ws.Range("D27").Value = "0.0₃0802"
$ws.Range("E27").Value = "  -2.52%  "

$ws.Range("D28").Value = "'7.16"
$ws.Range("E28").Value = "  -2.88%  "

$ws.Range("D29").Value = "'6.67"
$ws.Range("E29").Value = "  -3.26%  "

$ws.Range("D31").Value = "'1.59"
$ws.Range("E31").Value = "  -0.46%  "

$ws.Range("D32").Value = "'18.89"
$ws.Range("E32").Value = "  -1.49%  "

$ws.Range("D33").Value = "'150.64"
$ws.Range("E33").Value = "  +1.95%  "

$ws.Range("E34").Value = "  -3.83%  "

$ws.Range("E35").Value = "  -5.03%  "

$ws.Range("D36").Value = "'0.913"
$ws.Range("E36").Value = "  -4.92%  "

$ws.Range("D37").Value = "'0.867"
$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").Value = "'36.88"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("E39").Value = "  -4.80%  "

$ws.Range("E40").Value = "  -3.24%  "

$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").Value = "'0.613"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").Value = "'274.91"
$ws.Range("E43").Value = "  -4.55%  "

$ws.Range("D44").Value = "'19.79"
$ws.Range("E44").Value = "  -1.90%  "

$ws.Range("D45").Value = "'0.0967"
$ws.Range("E45").Value = "  -2.48%  "

$ws.Range("E46").Value = "  +2.02%  "

$ws.Range("D47").Value = "2.058.60"
$ws.Range("E47").Value = "  -3.99%  "

$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("E49").Value = "  -4.23%  "

$ws.Range("E50").Value = "  -2.94%  "

$ws.Range("D51").Value = "'18.77"
$ws.Range("E51").Value = "  -2.94%  "
